$wb = $excel.ActiveWorkbook

# --- Data edit -----------------------------------------------------------
# Hoja2!C9 and Hoja2!D9 ("JOSE" / "CAMILA" effort on sprint day 6) were
# blank; fill them in with 4 each. Every other changed cell in the diff
# (H9:J13, H14, J14 on Hoja2; L10:L15 on Hoja1; the chart's cached series)
# is a formula that depends on these two cells, so setting them and letting
# the workbook recalculate reproduces the rest of the diff automatically.
$wsHoja2 = $wb.Worksheets.Item("Hoja2")
$wsHoja2.Range("C9").Value = 4
$wsHoja2.Range("D9").Value = 4

# --- View / selection state ------------------------------------------------
# The saved workbook now opens on Hoja2 (activeTab moved from Hoja1 to
# Hoja2, tabSelected moved the same way) with F16 selected there, while
# Hoja1 keeps its previous selection (G12).
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
[void]$wsHoja1.Range("G12").Select()

[void]$wsHoja2.Activate()
[void]$wsHoja2.Range("F16").Select()
